# Scheduled market-data refresh: updates computed price/profit columns
# (H: currentAveragePrice, I: currentAveragePriceNQ, J: currentAveragePriceHQ,
#  K: LevePriceNQ, L: LevePriceHQ, M: LeveProfitNQ, N: LeveProfitHQ)
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1248.7778
$ws.Range("I28").Value = 1491.3334
$ws.Range("K28").Value = 1491.3334
$ws.Range("M28").Value = -1006.3334
$ws.Range("H64").Value = 7658
$ws.Range("J64").Value = 8865.200000000001
$ws.Range("L64").Value = 8865.200000000001
$ws.Range("N64").Value = -9361.200000000001
$ws.Range("H67").Value = 7658
$ws.Range("J67").Value = 8865.200000000001
$ws.Range("L67").Value = 8865.200000000001
$ws.Range("N67").Value = -10581.2
$ws.Range("H74").Value = 52169.523
$ws.Range("I74").Value = 67700
$ws.Range("J74").Value = 16671.285
$ws.Range("K74").Value = 67700
$ws.Range("L74").Value = 16671.285
$ws.Range("M74").Value = -66764
$ws.Range("N74").Value = -18543.285
$ws.Range("H77").Value = 52169.523
$ws.Range("I77").Value = 67700
$ws.Range("J77").Value = 16671.285
$ws.Range("K77").Value = 338500
$ws.Range("L77").Value = 83356.425
$ws.Range("M77").Value = -333820
$ws.Range("N77").Value = -92716.425
$ws.Range("H125").Value = 1182.7778
$ws.Range("I125").Value = 1122.1428
$ws.Range("K125").Value = 10099.2852
$ws.Range("M125").Value = -7639.2852
$ws.Range("H138").Value = 3082.02
$ws.Range("I138").Value = 1055.6666
$ws.Range("J138").Value = 3831.4932
$ws.Range("K138").Value = 3166.9998
$ws.Range("L138").Value = 11494.4796
$ws.Range("M138").Value = 1973.0002
$ws.Range("N138").Value = -21774.4796

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2316.3086
$ws.Range("I32").Value = 1914.2084
$ws.Range("K32").Value = 1914.2084
$ws.Range("M32").Value = -1627.2084
$ws.Range("H45").Value = 19022.467
$ws.Range("I45").Value = 16120.818
$ws.Range("K45").Value = 16120.818
$ws.Range("M45").Value = -15743.818
$ws.Range("H56").Value = 9583.333000000001
$ws.Range("J56").Value = 9583.333000000001
$ws.Range("L56").Value = 9583.333000000001
$ws.Range("N56").Value = -11067.333
$ws.Range("H125").Value = 29999
$ws.Range("J125").Value = 29999
$ws.Range("L125").Value = 29999
$ws.Range("N125").Value = -39839

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3005.4666
$ws.Range("I20").Value = 2470.7896
$ws.Range("J20").Value = 3929
$ws.Range("K20").Value = 2470.7896
$ws.Range("L20").Value = 3929
$ws.Range("M20").Value = -2223.7896
$ws.Range("N20").Value = -4423
$ws.Range("H134").Value = 2237.9119
$ws.Range("I134").Value = 1728.7407
$ws.Range("J134").Value = 4201.857
$ws.Range("K134").Value = 5186.2221
$ws.Range("L134").Value = 12605.571
$ws.Range("M134").Value = -2651.2221
$ws.Range("N134").Value = -17675.571

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 195
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("N22").ClearContents()
$ws.Range("H31").Value = 48351.434
$ws.Range("J31").Value = 10692.857
$ws.Range("L31").Value = 10692.857
$ws.Range("N31").Value = -11282.857
$ws.Range("H34").Value = 48351.434
$ws.Range("J34").Value = 10692.857
$ws.Range("L34").Value = 10692.857
$ws.Range("N34").Value = -11096.857
$ws.Range("H37").Value = 41028.5
$ws.Range("J37").Value = 41028.5
$ws.Range("L37").Value = 41028.5
$ws.Range("N37").Value = -41242.5
$ws.Range("H55").Value = 6949.5
$ws.Range("I55").Value = 6949.5
$ws.Range("K55").Value = 6949.5
$ws.Range("M55").Value = -6634.5
$ws.Range("H86").Value = 6100
$ws.Range("I86").Value = 6125
$ws.Range("J86").Value = 6000
$ws.Range("K86").Value = 6125
$ws.Range("L86").Value = 6000
$ws.Range("M86").Value = -5002
$ws.Range("N86").Value = -8246
$ws.Range("H89").Value = 6100
$ws.Range("I89").Value = 6125
$ws.Range("J89").Value = 6000
$ws.Range("K89").Value = 30625
$ws.Range("L89").Value = 30000
$ws.Range("M89").Value = -25009
$ws.Range("N89").Value = -41232
$ws.Range("H132").Value = 4121.56
$ws.Range("I132").Value = 3592.7273
$ws.Range("K132").Value = 10778.1819
$ws.Range("M132").Value = -8248.1819

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 44.033333
$ws.Range("J2").Value = 38.22222
$ws.Range("L2").Value = 229.33332
$ws.Range("N2").Value = -455.33332
$ws.Range("H7").Value = 173.88889
$ws.Range("J7").Value = 121.25
$ws.Range("L7").Value = 363.75
$ws.Range("N7").Value = -587.75
$ws.Range("I22").Value = 1500
$ws.Range("J22").Value = 5193.8
$ws.Range("K22").Value = 4500
$ws.Range("L22").Value = 15581.4
$ws.Range("M22").Value = -4331
$ws.Range("N22").Value = -15919.4
$ws.Range("I27").Value = 1500
$ws.Range("J27").Value = 5193.8
$ws.Range("K27").Value = 4500
$ws.Range("L27").Value = 15581.4
$ws.Range("M27").Value = -4398
$ws.Range("N27").Value = -15785.4
$ws.Range("H34").Value = 595096.2
$ws.Range("I34").Value = 1010419.5
$ws.Range("J34").Value = 1777.1428
$ws.Range("K34").Value = 3031258.5
$ws.Range("L34").Value = 5331.428400000001
$ws.Range("M34").Value = -3031174.5
$ws.Range("N34").Value = -5499.428400000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2309.7058
$ws.Range("I80").Value = 1793.1818
$ws.Range("K80").Value = 1793.1818
$ws.Range("M80").Value = -795.1818000000001
$ws.Range("H83").Value = 2309.7058
$ws.Range("I83").Value = 1793.1818
$ws.Range("K83").Value = 8965.909
$ws.Range("M83").Value = -3973.909
$ws.Range("H116").Value = 248000
$ws.Range("J116").Value = 248000
$ws.Range("L116").Value = 248000
$ws.Range("N116").Value = -257178
$ws.Range("H132").Value = 438191.6
$ws.Range("I132").Value = 438191.6
$ws.Range("K132").Value = 1314574.8
$ws.Range("M132").Value = -1312044.8
$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3000.182
$ws.Range("I46").Value = 2373
$ws.Range("J46").Value = 3522.8333
$ws.Range("K46").Value = 2373
$ws.Range("L46").Value = 3522.8333
$ws.Range("M46").Value = -2185
$ws.Range("N46").Value = -3898.8333
$ws.Range("H56").Value = 7249
$ws.Range("I56").Value = 7249
$ws.Range("K56").Value = 7249
$ws.Range("M56").Value = -6558
$ws.Range("H93").Value = 465448.34
$ws.Range("I93").Value = 619659.5
$ws.Range("K93").Value = 619659.5
$ws.Range("M93").Value = -618411.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H40").Value = 20000
$ws.Range("I40").Value = 20000
$ws.Range("K40").Value = 20000
$ws.Range("M40").Value = -19851
$ws.Range("H41").Value = 13847.929
$ws.Range("J41").Value = 15412
$ws.Range("L41").Value = 15412
$ws.Range("N41").Value = -16192
$ws.Range("H58").Value = 35000
$ws.Range("I58").Value = 35000
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 35000
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = -34692
$ws.Range("N58").ClearContents()
$ws.Range("H107").Value = 1427.5
$ws.Range("I107").Value = 904.625
$ws.Range("J107").Value = 1950.375
$ws.Range("K107").Value = 2713.875
$ws.Range("L107").Value = 5851.125
$ws.Range("M107").Value = -793.875
$ws.Range("N107").Value = -9691.125
$ws.Range("H109").Value = 39933.332
$ws.Range("J109").Value = 39933.332
$ws.Range("L109").Value = 39933.332
$ws.Range("N109").Value = -42707.332
$ws.Range("H129").Value = 33090.668
$ws.Range("I129").Value = 0
$ws.Range("J129").Value = 33090.668
$ws.Range("K129").Value = 0
$ws.Range("L129").Value = 33090.668
$ws.Range("M129").ClearContents()
$ws.Range("N129").Value = -43090.668

Write-Host "Updated leve profit sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)."

